$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 1942
    $ws.Range("F4").Value = 1203
    $ws.Range("F5").Value = 1349
    $ws.Range("F7").Value = 6062
}
